$wb = $excel.ActiveWorkbook

$ws_cash_flow = $wb.Worksheets.Item("Cash flow")
$ws_cash_flow.Range("E3").Value = 4.2987816484375
$ws_cash_flow.Range("J3").Value = 5.307375045358165
$ws_cash_flow.Range("K3").Value = 0.5528560567243939
$ws_cash_flow.Range("N3").Value = -9.053300637071271
$ws_cash_flow.Range("O3").Value = -4.75451898863377
$ws_cash_flow.Range("Q3").Value = -4.528113322508353
$ws_cash_flow.Range("R3").Value = -38.91836651000835
$ws_cash_flow.Range("E4").Value = 4.2987816484375
$ws_cash_flow.Range("J4").Value = 5.307375045358165
$ws_cash_flow.Range("K4").Value = 0.5528560567243939
$ws_cash_flow.Range("N4").Value = -9.053300637071271
$ws_cash_flow.Range("O4").Value = -4.75451898863377
$ws_cash_flow.Range("Q4").Value = -4.312488878579384
$ws_cash_flow.Range("R4").Value = -43.23085538858773
$ws_cash_flow.Range("E5").Value = 4.2987816484375
$ws_cash_flow.Range("J5").Value = 5.307375045358165
$ws_cash_flow.Range("K5").Value = 0.5528560567243939
$ws_cash_flow.Range("N5").Value = -9.053300637071271
$ws_cash_flow.Range("O5").Value = -4.75451898863377
$ws_cash_flow.Range("Q5").Value = -4.107132265313698
$ws_cash_flow.Range("R5").Value = -47.33798765390144
$ws_cash_flow.Range("E6").Value = 4.2987816484375
$ws_cash_flow.Range("J6").Value = 5.307375045358165
$ws_cash_flow.Range("K6").Value = 0.5528560567243939
$ws_cash_flow.Range("N6").Value = -9.053300637071271
$ws_cash_flow.Range("O6").Value = -4.75451898863377
$ws_cash_flow.Range("Q6").Value = -3.911554538393998
$ws_cash_flow.Range("R6").Value = -51.24954219229543
$ws_cash_flow.Range("E7").Value = 4.2987816484375
$ws_cash_flow.Range("J7").Value = 5.307375045358165
$ws_cash_flow.Range("K7").Value = 0.5528560567243939
$ws_cash_flow.Range("N7").Value = -9.053300637071271
$ws_cash_flow.Range("O7").Value = -4.75451898863377
$ws_cash_flow.Range("Q7").Value = -3.725290036565712
$ws_cash_flow.Range("R7").Value = -54.97483222886115
$ws_cash_flow.Range("E8").Value = 4.2987816484375
$ws_cash_flow.Range("J8").Value = 5.307375045358165
$ws_cash_flow.Range("K8").Value = 0.5528560567243939
$ws_cash_flow.Range("N8").Value = -9.053300637071271
$ws_cash_flow.Range("O8").Value = -4.75451898863377
$ws_cash_flow.Range("Q8").Value = -3.547895272919726
$ws_cash_flow.Range("R8").Value = -58.52272750178087
$ws_cash_flow.Range("E9").Value = 4.2987816484375
$ws_cash_flow.Range("J9").Value = 5.307375045358165
$ws_cash_flow.Range("K9").Value = 0.5528560567243939
$ws_cash_flow.Range("N9").Value = -9.053300637071271
$ws_cash_flow.Range("O9").Value = -4.75451898863377
$ws_cash_flow.Range("Q9").Value = -3.378947878971168
$ws_cash_flow.Range("R9").Value = -61.90167538075203
$ws_cash_flow.Range("E10").Value = 4.2987816484375
$ws_cash_flow.Range("J10").Value = 5.307375045358165
$ws_cash_flow.Range("K10").Value = 0.5528560567243939
$ws_cash_flow.Range("N10").Value = -9.053300637071271
$ws_cash_flow.Range("O10").Value = -4.75451898863377
$ws_cash_flow.Range("Q10").Value = -3.218045599020159
$ws_cash_flow.Range("R10").Value = -65.1197209797722

$ws_stream_table = $wb.Worksheets.Item("Stream table")
$ws_stream_table.Range("G7").Value = 0.217967709139475
$ws_stream_table.Range("H7").Value = 11.62381134201204
$ws_stream_table.Range("I7").Value = 0.05540882597631416
$ws_stream_table.Range("J7").Value = 26.03929351319754
$ws_stream_table.Range("K7").Value = 0.2733765351157892
$ws_stream_table.Range("M7").Value = 93.69229730444673
$ws_stream_table.Range("Q7").Value = 8.071623542867124
$ws_stream_table.Range("S7").Value = 7206.591606544994
$ws_stream_table.Range("X7").Value = 7233.005660068986
$ws_stream_table.Range("Y7").Value = 0
$ws_stream_table.Range("Z7").Value = 25040.08415496899
$ws_stream_table.Range("AD7").Value = 25034.91626845585
$ws_stream_table.Range("AF7").Value = 152.9049827971173
$ws_stream_table.Range("AG7").Value = 5.167886513142093
$ws_stream_table.Range("AH7").Value = 7645.249139855849
$ws_stream_table.Range("AL7").Value = 7492.344157058732
$ws_stream_table.Range("AO7").Value = 7492.344157058732
$ws_stream_table.Range("M9").Value = 82.36807616774904
$ws_stream_table.Range("Z9").Value = 0.3244168633124911
$ws_stream_table.Range("AD9").Value = 0.3082596401741453
$ws_stream_table.Range("AF9").Value = 0.03822699293088688
$ws_stream_table.Range("AG9").Value = 78.59523944631658
$ws_stream_table.Range("AH9").Value = 0.03822699293088698
$ws_stream_table.Range("AL9").Value = 0.03822699293088699
$ws_stream_table.Range("AO9").Value = 0.03822699293088699
$ws_stream_table.Range("M10").Value = 17.63192383225095
$ws_stream_table.Range("X10").Value = 0.1115943208426874
$ws_stream_table.Range("Z10").Value = 0.06731963135573926
$ws_stream_table.Range("AD10").Value = 0.06598685739652804
$ws_stream_table.Range("AF10").Value = 0.1636278020814877
$ws_stream_table.Range("AG10").Value = 6.523708406298974
$ws_stream_table.Range("AH10").Value = 0.1636278020814881
$ws_stream_table.Range("AL10").Value = 0.1636278020814881
$ws_stream_table.Range("AO10").Value = 0.1636278020814881
$ws_stream_table.Range("X11").Value = 0.1028243047317388
$ws_stream_table.Range("Z11").Value = 0.01471494561359183
$ws_stream_table.Range("AD11").Value = 0.01442362351437578
$ws_stream_table.Range("AF11").Value = 0.09926516488389507
$ws_stream_table.Range("AG11").Value = 1.425973560228607
$ws_stream_table.Range("AH11").Value = 0.09926516488389536
$ws_stream_table.Range("AL11").Value = 0.09926516488389538
$ws_stream_table.Range("AO11").Value = 0.09926516488389538
$ws_stream_table.Range("X12").Value = 0.15076916447274
$ws_stream_table.Range("Z12").Value = 0.1239521090390938
$ws_stream_table.Range("AD12").Value = 0.1214981421977776
$ws_stream_table.Range("AF12").Value = 0.1455504708720297
$ws_stream_table.Range("AG12").Value = 12.01176238538453
$ws_stream_table.Range("AH12").Value = 0.1455504708720299
$ws_stream_table.Range("AL12").Value = 0.1455504708720299
$ws_stream_table.Range("AO12").Value = 0.1455504708720299
$ws_stream_table.Range("S13").Value = 0.06468812460663789
$ws_stream_table.Range("X13").Value = 0.06445189147396396
$ws_stream_table.Range("Z13").Value = 0.01489390827756699
$ws_stream_table.Range("AD13").Value = 0.01459904312896975
$ws_stream_table.Range("AF13").Value = 0.06222096663754162
$ws_stream_table.Range("AG13").Value = 1.443316201771296
$ws_stream_table.Range("AH13").Value = 0.06222096663754189
$ws_stream_table.Range("AL13").Value = 0.0622209666375419
$ws_stream_table.Range("AO13").Value = 0.0622209666375419
$ws_stream_table.Range("S14").Value = 0.4916297470104479
$ws_stream_table.Range("X14").Value = 0.489834375202126
$ws_stream_table.Range("AF14").Value = 0.4728793464453172
$ws_stream_table.Range("AH14").Value = 0.4728793464453184
$ws_stream_table.Range("AL14").Value = 0.4728793464453184
$ws_stream_table.Range("AO14").Value = 0.4728793464453184
$ws_stream_table.Range("S15").Value = 53.30093987678187
$ws_stream_table.Range("X15").Value = 53.10629135237148
$ws_stream_table.Range("Z15").Value = 99.04449004582048
$ws_stream_table.Range("AD15").Value = 99.06493551800902
$ws_stream_table.Range("AF15").Value = 51.2680808415731
$ws_stream_table.Range("AH15").Value = 51.26808084157318
$ws_stream_table.Range("AL15").Value = 51.26808084157319
$ws_stream_table.Range("AO15").Value = 51.26808084157319
$ws_stream_table.Range("S16").Value = 5.079447328322908
$ws_stream_table.Range("X16").Value = 5.060897806877995
$ws_stream_table.Range("Z16").Value = 0.410212496581023
$ws_stream_table.Range("AD16").Value = 0.4102971755791819
$ws_stream_table.Range("AF16").Value = 8.252977658961324
$ws_stream_table.Range("AH16").Value = 8.252977658961338
$ws_stream_table.Range("AL16").Value = 8.252977658961338
$ws_stream_table.Range("AO16").Value = 8.252977658961338
$ws_stream_table.Range("S23").Value = 1.748647314346649
$ws_stream_table.Range("X23").Value = 1.742261467863667
$ws_stream_table.Range("AF23").Value = 1.681955179075094
$ws_stream_table.Range("AH23").Value = 1.6819551790751
$ws_stream_table.Range("AL23").Value = 1.681955179075101
$ws_stream_table.Range("AO23").Value = 1.681955179075101
$ws_stream_table.Range("S24").Value = 39.3146476089315
$ws_stream_table.Range("X24").Value = 39.17107531616362
$ws_stream_table.Range("AF24").Value = 37.81521557653933
$ws_stream_table.Range("AH24").Value = 37.81521557653922
$ws_stream_table.Range("AL24").Value = 37.81521557653922
$ws_stream_table.Range("AO24").Value = 37.81521557653922

$ws_design_requirements = $wb.Worksheets.Item("Design requirements")
$ws_design_requirements.Range("E31").Value = 54.61668379141244
$ws_design_requirements.Range("D39").Value = 0.001616139599996895
